$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that needs to move from
# 2024-09-01 (45536) to 2024-09-02 (45537) for every data row (rows 2-28).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45536) {
        $cell.Value = 45537
    }
}
